$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.021506968033493
$ws.Range("D2").Value = 1.029927289697108
$ws.Range("E2").Value = 1.022360159033294
$ws.Range("F2").Value = 1.038399514175351
$ws.Range("I2").Value = 1.026523612493456
$ws.Range("J2").Value = 1.026698047430847
$ws.Range("K2").Value = 1.032739781674504
$ws.Range("L2").Value = 1.025194769800136
$ws.Range("M2").Value = 1.041187657277945
$ws.Range("N2").Value = 1.012906584744852
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.022549443408616
$ws.Range("D3").Value = 1.03086100952761
$ws.Range("E3").Value = 1.023248168031278
$ws.Range("F3").Value = 1.039420861456714
$ws.Range("I3").Value = 1.026487278619574
$ws.Range("J3").Value = 1.027377687354206
$ws.Range("K3").Value = 1.033481411957718
$ws.Range("L3").Value = 1.02588920599141
$ws.Range("M3").Value = 1.042018451621703
$ws.Range("N3").Value = 1.013137385300332
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.023224573598764
$ws.Range("D4").Value = 1.031465864900698
$ws.Range("E4").Value = 1.023823619344901
$ws.Range("F4").Value = 1.04008196648891
$ws.Range("I4").Value = 1.026461399693601
$ws.Range("J4").Value = 1.027817494128044
$ws.Range("K4").Value = 1.033961350299651
$ws.Range("L4").Value = 1.026338771949655
$ws.Range("M4").Value = 1.042555651953117
$ws.Range("N4").Value = 1.013286587367788
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023508536633653
$ws.Range("D5").Value = 1.031720307056358
$ws.Range("E5").Value = 1.024065741701133
$ws.Range("F5").Value = 1.040359947728098
$ws.Range("I5").Value = 1.026449952123637
$ws.Range("J5").Value = 1.028002396493076
$ws.Range("K5").Value = 1.034163128515313
$ws.Range("L5").Value = 1.026527821297644
$ws.Range("M5").Value = 1.04278139901955
$ws.Range("N5").Value = 1.013349277772948
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.023556223378539
$ws.Range("D6").Value = 1.031763038425489
$ws.Range("E6").Value = 1.02410640696611
$ws.Range("F6").Value = 1.040406625056001
$ws.Range("I6").Value = 1.026447996683618
$ws.Range("J6").Value = 1.028033442841466
$ws.Range("K6").Value = 1.034197008654589
$ws.Range("L6").Value = 1.026559566533876
$ws.Range("M6").Value = 1.042819297518781
$ws.Range("N6").Value = 1.013359801756575
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023228367383658
$ws.Range("D7").Value = 1.031469264137963
$ws.Range("E7").Value = 1.02382685380063
$ws.Range("F7").Value = 1.040085680682352
$ws.Range("I7").Value = 1.026461248964081
$ws.Range("J7").Value = 1.02781996477342
$ws.Range("K7").Value = 1.033964046423349
$ws.Range("L7").Value = 1.026341297833546
$ws.Range("M7").Value = 1.042558668759297
$ws.Range("N7").Value = 1.013287425174604
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.021859157423555
$ws.Range("D8").Value = 1.03024270411428
$ws.Range("E8").Value = 1.022660089272415
$ws.Range("F8").Value = 1.038744635957469
$ws.Range("I8").Value = 1.026511823072012
$ws.Range("J8").Value = 1.026927727454551
$ws.Range("K8").Value = 1.032990407459508
$ws.Range("L8").Value = 1.025429411942423
$ws.Range("M8").Value = 1.041468505915464
$ws.Range("N8").Value = 1.012984613948929
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019450880926046
$ws.Range("D9").Value = 1.028086568338479
$ws.Range("E9").Value = 1.020610651281975
$ws.Range("F9").Value = 1.036383319683319
$ws.Range("I9").Value = 1.02658284726232
$ws.Range("J9").Value = 1.025355780877988
$ws.Range("K9").Value = 1.031275178701595
$ws.Range("L9").Value = 1.023824262127868
$ws.Range("M9").Value = 1.039544647529762
$ws.Range("N9").Value = 1.012449953130937
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.017848367401393
$ws.Range("D10").Value = 1.02665270728313
$ws.Range("E10").Value = 1.019248814922371
$ws.Range("F10").Value = 1.034810371398216
$ws.Range("I10").Value = 1.026618092753561
$ws.Range("J10").Value = 1.024308046089744
$ws.Range("K10").Value = 1.03013203846829
$ws.Range("L10").Value = 1.022755350899078
$ws.Range("M10").Value = 1.038260225766884
$ws.Range("N10").Value = 1.012092811949419
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.01715517620961
$ws.Range("D11").Value = 1.026032685405129
$ws.Range("E11").Value = 1.0186601920162
$ws.Range("F11").Value = 1.034129580727667
$ws.Range("I11").Value = 1.026630496300024
$ws.Range("J11").Value = 1.023854426729052
$ws.Range("K11").Value = 1.029637137066868
$ws.Range("L11").Value = 1.022292790319068
$ws.Range("M11").Value = 1.037703632116647
$ws.Range("N11").Value = 1.011938003327347
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016897800559698
$ws.Range("D12").Value = 1.025802510098095
$ws.Range("E12").Value = 1.018441711493375
$ws.Range("F12").Value = 1.033876751823956
$ws.Range("I12").Value = 1.02663467514352
$ws.Range("J12").Value = 1.023685941323487
$ws.Range("K12").Value = 1.029453322216421
$ws.Range("L12").Value = 1.022121018130806
$ws.Range("M12").Value = 1.037496824814314
$ws.Range("N12").Value = 1.011880476126287
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.016953003692869
$ws.Range("D13").Value = 1.025851877658316
$ws.Range("E13").Value = 1.018488569043507
$ws.Range("F13").Value = 1.033930982335604
$ws.Range("I13").Value = 1.026633798136936
$ws.Range("J13").Value = 1.023722081598577
$ws.Range("K13").Value = 1.029492750501263
$ws.Range("L13").Value = 1.022157861874669
$ws.Range("M13").Value = 1.037541188545674
$ws.Range("N13").Value = 1.011892817005328
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.017133899286544
$ws.Range("D14").Value = 1.026013656417939
$ws.Range("E14").Value = 1.018642129074719
$ws.Range("F14").Value = 1.034108680845515
$ws.Range("I14").Value = 1.026630850457769
$ws.Range("J14").Value = 1.023840499485912
$ws.Range("K14").Value = 1.029621942593409
$ws.Range("L14").Value = 1.022278590687544
$ws.Range("M14").Value = 1.037686538652683
$ws.Range("N14").Value = 1.01193324860685
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.017245369137085
$ws.Range("D15").Value = 1.026113350692862
$ws.Range("E15").Value = 1.018736763734309
$ws.Range("F15").Value = 1.034218173029785
$ws.Range("I15").Value = 1.026628977559578
$ws.Range("J15").Value = 1.02391346186515
$ws.Range("K15").Value = 1.029701543922724
$ws.Range("L15").Value = 1.022352981462971
$ws.Range("M15").Value = 1.037776085249338
$ws.Range("N15").Value = 1.01195815662775
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017894387136913
$ws.Range("D16").Value = 1.026693874058568
$ws.Range("E16").Value = 1.019287902312787
$ws.Range("F16").Value = 1.034855559775713
$ws.Range("I16").Value = 1.0266172094417
$ws.Range("J16").Value = 1.024338152550487
$ws.Range("K16").Value = 1.030164885310241
$ws.Range("L16").Value = 1.022786055576518
$ws.Range("M16").Value = 1.038297156072597
$ws.Range("N16").Value = 1.012103082653041
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.018301688527675
$ws.Range("D17").Value = 1.027058249014123
$ws.Range("E17").Value = 1.019633901664018
$ws.Range("F17").Value = 1.035255458388853
$ws.Range("I17").Value = 1.026609063154189
$ws.Range("J17").Value = 1.024604565335694
$ws.Range("K17").Value = 1.030455550579408
$ws.Range("L17").Value = 1.023057788250155
$ws.Range("M17").Value = 1.038623895503254
$ws.Range("N17").Value = 1.012193947169606
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018539328817957
$ws.Range("D18").Value = 1.027270864663507
$ws.Range("E18").Value = 1.019835819571095
$ws.Range("F18").Value = 1.035488741969415
$ws.Range("I18").Value = 1.026604035742517
$ws.Range("J18").Value = 1.024759964844608
$ws.Range("K18").Value = 1.030625098953152
$ws.Range("L18").Value = 1.023216312784
$ws.Range("M18").Value = 1.038814435473126
$ws.Range("N18").Value = 1.012246931049981
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018620369594491
$ws.Range("D19").Value = 1.027343374972949
$ws.Range("E19").Value = 1.019904685693701
$ws.Range("F19").Value = 1.035568290611265
$ws.Range("I19").Value = 1.02660227471079
$ws.Range("J19").Value = 1.024812952968517
$ws.Range("K19").Value = 1.030682911944464
$ws.Range("L19").Value = 1.023270370193812
$ws.Range("M19").Value = 1.03887939758885
$ws.Range("N19").Value = 1.01226499448639
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.01825798187981
$ws.Range("D20").Value = 1.02701914655232
$ws.Range("E20").Value = 1.019596768616889
$ws.Range("F20").Value = 1.03521254999483
$ws.Range("I20").Value = 1.026609965695885
$ws.Range("J20").Value = 1.024575981198918
$ws.Range("K20").Value = 1.030424364089686
$ws.Range("L20").Value = 1.023028631061064
$ws.Range("M20").Value = 1.038588843751046
$ws.Range("N20").Value = 1.012184199904858
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.017080627104338
$ws.Range("D21").Value = 1.025966013053807
$ws.Range("E21").Value = 1.018596905043303
$ws.Range("F21").Value = 1.034056351755817
$ws.Range("I21").Value = 1.02663173029344
$ws.Range("J21").Value = 1.023805628116488
$ws.Range("K21").Value = 1.029583898362481
$ws.Range("L21").Value = 1.022243037858316
$ws.Range("M21").Value = 1.037643738417261
$ws.Range("N21").Value = 1.011921343178912
$ws.Range("B22").Value = 1.019999999999999
$ws.Range("C22").Value = 1.016340993208924
$ws.Range("D22").Value = 1.025304609651946
$ws.Range("E22").Value = 1.017969178582918
$ws.Range("F22").Value = 1.033329676843075
$ws.Range("I22").Value = 1.026642936534324
$ws.Range("J22").Value = 1.023321329007943
$ws.Range("K22").Value = 1.029055543017028
$ws.Range("L22").Value = 1.021749355658358
$ws.Range("M22").Value = 1.037049145507087
$ws.Range("N22").Value = 1.011755933767917
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.016733028729847
$ws.Range("D23").Value = 1.025655161371073
$ws.Range("E23").Value = 1.018301860078893
$ws.Range("F23").Value = 1.033714874756802
$ws.Range("I23").Value = 1.026637230453435
$ws.Range("J23").Value = 1.023578059989432
$ws.Range("K23").Value = 1.029335626489585
$ws.Range("L23").Value = 1.022011041906123
$ws.Range("M23").Value = 1.037364384961577
$ws.Range("N23").Value = 1.011843633696053
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.018277730821552
$ws.Range("D24").Value = 1.027036815018599
$ws.Range("E24").Value = 1.019613547126323
$ws.Range("F24").Value = 1.035231938360003
$ws.Range("I24").Value = 1.02660955872848
$ws.Range("J24").Value = 1.02458889712311
$ws.Range("K24").Value = 1.030438455896512
$ws.Range("L24").Value = 1.023041805854899
$ws.Range("M24").Value = 1.038604682257658
$ws.Range("N24").Value = 1.012188604323229
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.020072949714488
$ws.Range("D25").Value = 1.028643357056416
$ws.Range("E25").Value = 1.021139698573623
$ws.Range("F25").Value = 1.03699355901353
$ws.Range("I25").Value = 1.026566623658166
$ws.Range("J25").Value = 1.025762128835807
$ws.Range("K25").Value = 1.031718548806392
$ws.Range("L25").Value = 1.024239025192359
$ws.Range("M25").Value = 1.040042341998809
$ws.Range("N25").Value = 1.012588300571156
